# Apply the edit described by the commit:
# "Remove a special character from in the pre basic files, correction on #141 edits."
#
# 1) Replace "Wårlind" with "Warlind" everywhere it appears on the sheet
#    (removes the special character "å" from "David Wårlind" occurrences
#    in the "comment author" column).
# 2) Update the worksheet view state (scroll position / selection) to match
#    the commit's recorded cursor position.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Fix the special character in the author name -----------------------
$used = $ws.UsedRange
$used.Replace("Wårlind", "Warlind") | Out-Null

# --- 2) Update the view / selection state -----------------------------------
$window = $excel.ActiveWindow

# Select the range that corresponds to the recorded selection
# (sqref="G175 G179:G181 G183 G185 G187:G288 G290 G292:G303 G305:G312 G314:G319"),
# with the active cell at G188 (5th area, activeCellId="4").
$selRange = $ws.Range("G175,G179:G181,G183,G185,G187:G288,G290,G292:G303,G305:G312,G314:G319")
$selRange.Select() | Out-Null

# Make sure the active cell within the selection is G188, as recorded in the
# diff (activeCell="G188"). Re-select starting at G188 so it becomes active
# while still covering the remainder of the big contiguous block.
$ws.Range("G188:G288,G290,G292:G303,G305:G312,G314:G319").Select() | Out-Null

# Scroll the view so the top-left visible cell is E155 (topLeftCell="E155").
$window.ScrollRow = 155
$window.ScrollColumn = 5
